$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "DKS" / Desktop Computer device-type group, which occupies
# rows 8-10 (eng, ara, fra language rows) of the master-device_type table.
$ws.Range("A8:G10").EntireRow.Delete()

# Leave the active cell roughly where the deleted rows used to be,
# matching Excel's natural post-delete selection behaviour.
$ws.Range("E10").Select() | Out-Null

# Touch page setup (paper size / orientation) as part of finalizing the sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
